# Scheduled-runner style refresh of market-price / profit columns
# (H: currentAveragePrice, I: currentAveragePriceNQ, J: currentAveragePriceHQ,
#  K: LevePriceNQ, L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7115.2666
$ws.Range("I62").Value = 6089.125
$ws.Range("K62").Value = 6089.125
$ws.Range("M62").Value = -5465.125

$ws.Range("H65").Value = 7115.2666
$ws.Range("I65").Value = 6089.125
$ws.Range("K65").Value = 30445.625
$ws.Range("M65").Value = -27325.625

$ws.Range("H100").Value = 124457.2
$ws.Range("I100").Value = 201201.67
$ws.Range("K100").Value = 201201.67
$ws.Range("M100").Value = -200660.67

$ws.Range("H101").Value = 391.33334
$ws.Range("I101").Value = 337.25
$ws.Range("K101").Value = 1011.75
$ws.Range("M101").Value = 610.25

$ws.Range("H113").Value = 6186.7144
$ws.Range("J113").Value = 6741.4
$ws.Range("L113").Value = 6741.4
$ws.Range("N113").Value = -13249.4

$ws.Range("H131").Value = 4953.0713
$ws.Range("I131").Value = 4719.273
$ws.Range("J131").Value = 5810.3335
$ws.Range("K131").Value = 14157.819
$ws.Range("L131").Value = 17431.0005
$ws.Range("M131").Value = -9117.819
$ws.Range("N131").Value = -27511.0005

$ws.Range("H132").Value = 24292.193
$ws.Range("I132").Value = 24935.268
$ws.Range("K132").Value = 74805.804
$ws.Range("M132").Value = -72275.804

$ws.Range("H137").Value = 19495.666
$ws.Range("I137").Value = 29447.818
$ws.Range("J137").Value = 3856.5715
$ws.Range("K137").Value = 88343.454
$ws.Range("L137").Value = 11569.7145
$ws.Range("M137").Value = -85793.454
$ws.Range("N137").Value = -16669.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 16874.5
$ws.Range("I25").Value = 249
$ws.Range("J25").Value = 20199.6
$ws.Range("K25").Value = 249
$ws.Range("L25").Value = 20199.6
$ws.Range("M25").Value = 153
$ws.Range("N25").Value = -21003.6

$ws.Range("H37").Value = 22856.857
$ws.Range("I37").Value = 18333.166
$ws.Range("K37").Value = 18333.166
$ws.Range("M37").Value = -18060.166

$ws.Range("H44").Value = 59999
$ws.Range("J44").Value = 59999
$ws.Range("L44").Value = 59999
$ws.Range("N44").Value = -60975

$ws.Range("H97").Value = 1578.7241
$ws.Range("I97").Value = 1067.7273
$ws.Range("K97").Value = 1067.7273
$ws.Range("M97").Value = -571.7273

$ws.Range("H102").Value = 4734.0527
$ws.Range("I102").Value = 4163.722
$ws.Range("K102").Value = 4163.722
$ws.Range("M102").Value = -2541.722

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 497.07144
$ws.Range("I22").Value = 504.53845
$ws.Range("K22").Value = 504.53845
$ws.Range("M22").Value = -331.53845

$ws.Range("H24").Value = 300
$ws.Range("I24").Value = 300
$ws.Range("K24").Value = 300
$ws.Range("M24").Value = -65

$ws.Range("H80").Value = 1958.5625
$ws.Range("I80").Value = 1593.5555
$ws.Range("J80").Value = 2427.8572
$ws.Range("K80").Value = 1593.5555
$ws.Range("L80").Value = 2427.8572
$ws.Range("M80").Value = -595.5554999999999
$ws.Range("N80").Value = -4423.8572

$ws.Range("H83").Value = 1958.5625
$ws.Range("I83").Value = 1593.5555
$ws.Range("J83").Value = 2427.8572
$ws.Range("K83").Value = 7967.7775
$ws.Range("L83").Value = 12139.286
$ws.Range("M83").Value = -2975.7775
$ws.Range("N83").Value = -22123.286

$ws.Range("H105").Value = 3804.2144
$ws.Range("I105").Value = 3105
$ws.Range("K105").Value = 3105
$ws.Range("M105").Value = -1358

$ws.Range("H134").Value = 1983.1714
$ws.Range("I134").Value = 1771.3226
$ws.Range("K134").Value = 5313.9678
$ws.Range("M134").Value = -2778.9678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4001372.5
$ws.Range("I31").Value = 4168013
$ws.Range("K31").Value = 4168013
$ws.Range("M31").Value = -4167718

$ws.Range("H34").Value = 4001372.5
$ws.Range("I34").Value = 4168013
$ws.Range("K34").Value = 4168013
$ws.Range("M34").Value = -4167811

$ws.Range("H41").Value = 11708.272
$ws.Range("I41").Value = 1231.875
$ws.Range("K41").Value = 1231.875
$ws.Range("M41").Value = -803.875

$ws.Range("H50").Value = 38152.5
$ws.Range("J50").Value = 38152.5
$ws.Range("L50").Value = 38152.5
$ws.Range("N50").Value = -39402.5

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H59").Value = 24998
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 34199.5
$ws.Range("I60").Value = 34499
$ws.Range("J60").Value = 33900
$ws.Range("K60").Value = 34499
$ws.Range("L60").Value = 33900
$ws.Range("M60").Value = -33988
$ws.Range("N60").Value = -34922

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 74998.5

$ws.Range("H77").Value = 74998.5

$ws.Range("H132").Value = 41271.4
$ws.Range("I132").Value = 47144
$ws.Range("J132").Value = 3099.5
$ws.Range("K132").Value = 141432
$ws.Range("L132").Value = 9298.5
$ws.Range("M132").Value = -138902
$ws.Range("N132").Value = -14358.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79532200
$ws.Range("J4").Value = 96364910
$ws.Range("L4").Value = 289094730
$ws.Range("N4").Value = -289094954

$ws.Range("H24").Value = 3567.3333
$ws.Range("I24").Value = 5100
$ws.Range("J24").Value = 502
$ws.Range("K24").Value = 15300
$ws.Range("L24").Value = 1506
$ws.Range("M24").Value = -15070
$ws.Range("N24").Value = -1966

$ws.Range("H96").Value = 166670000

$ws.Range("H114").Value = 834.3
$ws.Range("J114").Value = 1387.5
$ws.Range("L114").Value = 4162.5
$ws.Range("N114").Value = -10670.5

$ws.Range("H137").Value = 2493
$ws.Range("J137").Value = 1433
$ws.Range("L137").Value = 4299
$ws.Range("N137").Value = -14499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6362.3
$ws.Range("I43").Value = 883.8570999999999
$ws.Range("K43").Value = 883.8570999999999
$ws.Range("M43").Value = -732.8570999999999

$ws.Range("H46").Value = 34999.5
$ws.Range("I46").Value = 30000
$ws.Range("J46").Value = 39999
$ws.Range("K46").Value = 30000
$ws.Range("L46").Value = 39999
$ws.Range("M46").Value = -29844
$ws.Range("N46").Value = -40311

$ws.Range("H122").Value = 4016
$ws.Range("I122").Value = 3597.111
$ws.Range("K122").Value = 10791.333
$ws.Range("M122").Value = -8341.332999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4717.55
$ws.Range("I46").Value = 1075.625
$ws.Range("J46").Value = 7145.5
$ws.Range("K46").Value = 1075.625
$ws.Range("L46").Value = 7145.5
$ws.Range("M46").Value = -887.625
$ws.Range("N46").Value = -7521.5

$ws.Range("H93").Value = 3039.4
$ws.Range("I93").Value = 2335.6875
$ws.Range("J93").Value = 4290.4443
$ws.Range("K93").Value = 2335.6875
$ws.Range("L93").Value = 4290.4443
$ws.Range("M93").Value = -1087.6875
$ws.Range("N93").Value = -6786.4443

$ws.Range("H100").Value = 3261.5
$ws.Range("J100").Value = 3492.25
$ws.Range("L100").Value = 3492.25
$ws.Range("N100").Value = -4574.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15430.286
$ws.Range("J54").Value = 37423
$ws.Range("L54").Value = 37423
$ws.Range("N54").Value = -38463

$ws.Range("H81").Value = 2735.3635
$ws.Range("J81").Value = 842.5
$ws.Range("L81").Value = 1685
$ws.Range("N81").Value = -3807

$ws.Range("H84").Value = 2735.3635
$ws.Range("J84").Value = 842.5
$ws.Range("L84").Value = 8425
$ws.Range("N84").Value = -19033

$ws.Range("H100").Value = 1166.5
$ws.Range("I100").Value = 624.75
$ws.Range("K100").Value = 1249.5
$ws.Range("M100").Value = -708.5

$ws.Range("H136").Value = 21440.97
$ws.Range("I136").Value = 24915.414
$ws.Range("K136").Value = 74746.242
$ws.Range("M136").Value = -72196.242
